$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Step 1: simple text replacements in the first three rows ---
$t.Cell(1, 1).Range.Text = "0M"
$t.Cell(2, 1).Range.Text = "0M"
$t.Cell(3, 1).Range.Text = "0M"

# --- Step 2: insert 10 new rows before (current) row 4, each holding a
# single value that used to live, tab-separated, inside one run. Insert
# them in reverse order so each new row lands directly above the next,
# preserving the intended top-to-bottom order.
$newValues = @("21", "0.00002", "0.00004", "0.00003", "0.00001", "0.00003", "0.00003", "0.00004", "0.00066", "100.0")

$refRow = $t.Rows.Item(4)
for ($i = $newValues.Length - 1; $i -ge 0; $i--) {
    $newRow = $t.Rows.Add($refRow)
    $newRow.Cells.Item(1).Range.Text = $newValues[$i]
    $refRow = $newRow
}

# --- Step 3: the two tab-separated multi-run rows collapse down to a
# single value each. After the inserts above, these are rows 44 and 45.
$t.Cell(44, 1).Range.Text = "100"
$t.Cell(45, 1).Range.Text = "0"

# --- Step 4: the trailing, previously-empty row gets new text.
$t.Cell(46, 1).Range.Text = "70"
